$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 3432.2222
$ws.Range("J80").Value = 4119.6
$ws.Range("L80").Value = 12358.8
$ws.Range("N80").Value = -14354.8
$ws.Range("H83").Value = 3432.2222
$ws.Range("J83").Value = 4119.6
$ws.Range("L83").Value = 37076.4
$ws.Range("N83").Value = -47060.4
$ws.Range("H98").Value = 941.2222
$ws.Range("I98").Value = 941.2222
$ws.Range("K98").Value = 941.2222
$ws.Range("M98").Value = 556.7778
$ws.Range("H107").Value = 647.8
$ws.Range("I107").Value = 322.875
$ws.Range("J107").Value = 1947.5
$ws.Range("K107").Value = 322.875
$ws.Range("L107").Value = 1947.5
$ws.Range("M107").Value = 1597.125
$ws.Range("N107").Value = -5787.5
$ws.Range("H122").Value = 941.2222
$ws.Range("I122").Value = 941.2222
$ws.Range("K122").Value = 2823.6666
$ws.Range("M122").Value = -373.6666
$ws.Range("H138").Value = 1447.2
$ws.Range("I138").Value = 670.06665
$ws.Range("J138").Value = 3778.6
$ws.Range("K138").Value = 2010.19995
$ws.Range("L138").Value = 11335.8
$ws.Range("M138").Value = 3129.80005
$ws.Range("N138").Value = -21615.8

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 303.66666
$ws.Range("I97").Value = 303.66666
$ws.Range("J97").Value = 0
$ws.Range("K97").Value = 303.66666
$ws.Range("L97").Value = 0
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = 192.33334
$ws.Range("H122").Value = 1519.2307
$ws.Range("I122").Value = 1519.2307
$ws.Range("K122").Value = 4557.6921
$ws.Range("M122").Value = -2107.6921
$ws.Range("H132").Value = 9673
$ws.Range("I132").Value = 9632.125
$ws.Range("K132").Value = 28896.375
$ws.Range("M132").Value = -26366.375

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 4705.8
$ws.Range("I134").Value = 4705.8
$ws.Range("K134").Value = 14117.4
$ws.Range("M134").Value = -11582.4

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1427.875
$ws.Range("I31").Value = 1704
$ws.Range("J31").Value = 599.5
$ws.Range("K31").Value = 1704
$ws.Range("L31").Value = 599.5
$ws.Range("M31").Value = -1409
$ws.Range("N31").Value = -1189.5
$ws.Range("H34").Value = 1427.875
$ws.Range("I34").Value = 1704
$ws.Range("J34").Value = 599.5
$ws.Range("K34").Value = 1704
$ws.Range("L34").Value = 599.5
$ws.Range("M34").Value = -1502
$ws.Range("N34").Value = -1003.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1073.6
$ws.Range("I5").Value = 993.375
$ws.Range("J5").Value = 1394.5
$ws.Range("K5").Value = 2980.125
$ws.Range("L5").Value = 4183.5
$ws.Range("M5").Value = -2868.125
$ws.Range("N5").Value = -4407.5
$ws.Range("H81").Value = 11239
$ws.Range("J81").Value = 13625
$ws.Range("L81").Value = 40875
$ws.Range("N81").Value = -43121
$ws.Range("H84").Value = 11239
$ws.Range("J84").Value = 13625
$ws.Range("L84").Value = 122625
$ws.Range("N84").Value = -133857
$ws.Range("H107").Value = 786.1
$ws.Range("J107").Value = 790.6
$ws.Range("L107").Value = 2371.8
$ws.Range("N107").Value = -6211.8
$ws.Range("H129").Value = 2797
$ws.Range("J129").Value = 2797
$ws.Range("L129").Value = 8391
$ws.Range("N129").Value = -18391
$ws.Range("H131").Value = 47629.766
$ws.Range("I131").Value = 1447.1111
$ws.Range("J131").Value = 99585.25
$ws.Range("K131").Value = 4341.3333
$ws.Range("L131").Value = 298755.75
$ws.Range("M131").Value = 698.6666999999998
$ws.Range("N131").Value = -308835.75
$ws.Range("H135").Value = 1073.6
$ws.Range("I135").Value = 993.375
$ws.Range("J135").Value = 1394.5
$ws.Range("K135").Value = 8940.375
$ws.Range("L135").Value = 12550.5
$ws.Range("M135").Value = -6405.375
$ws.Range("N135").Value = -17620.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H58").Value = 10000
$ws.Range("I58").Value = 10000
$ws.Range("K58").Value = 10000
$ws.Range("M58").Value = -9723

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5763.905
$ws.Range("I7").Value = 5714.8823
$ws.Range("J7").Value = 5972.25
$ws.Range("K7").Value = 5714.8823
$ws.Range("L7").Value = 5972.25
$ws.Range("M7").Value = -5602.8823
$ws.Range("N7").Value = -6196.25
$ws.Range("H61").Value = 1902
$ws.Range("I61").Value = 1890
$ws.Range("J61").Value = 1998
$ws.Range("K61").Value = 1890
$ws.Range("L61").Value = 1998
$ws.Range("M61").Value = -1688
$ws.Range("N61").Value = -2402
$ws.Range("H93").Value = 999.5
$ws.Range("I93").Value = 999
$ws.Range("K93").Value = 999
$ws.Range("M93").Value = 249
$ws.Range("H113").Value = 1902
$ws.Range("I113").Value = 1890
$ws.Range("J113").Value = 1998
$ws.Range("K113").Value = 1890
$ws.Range("L113").Value = 1998
$ws.Range("M113").Value = 280
$ws.Range("N113").Value = -6338
$ws.Range("H122").Value = 3818.2222
$ws.Range("I122").Value = 3548.6155
$ws.Range("J122").Value = 4519.2
$ws.Range("K122").Value = 10645.8465
$ws.Range("L122").Value = 13557.6
$ws.Range("M122").Value = -8195.8465
$ws.Range("N122").Value = -18457.6
$ws.Range("H126").Value = 5763.905
$ws.Range("I126").Value = 5714.8823
$ws.Range("J126").Value = 5972.25
$ws.Range("K126").Value = 17144.6469
$ws.Range("L126").Value = 17916.75
$ws.Range("M126").Value = -14674.6469
$ws.Range("N126").Value = -22856.75
$ws.Range("H132").Value = 7360.4
$ws.Range("I132").Value = 7201.5
$ws.Range("J132").Value = 7466.3335
$ws.Range("K132").Value = 21604.5
$ws.Range("L132").Value = 22399.0005
$ws.Range("M132").Value = -19074.5
$ws.Range("N132").Value = -27459.0005
$ws.Range("H136").Value = 2802.25
$ws.Range("I136").Value = 1295.0769
$ws.Range("J136").Value = 9333.333000000001
$ws.Range("K136").Value = 3885.2307
$ws.Range("L136").Value = 27999.999
$ws.Range("M136").Value = -1335.2307
$ws.Range("N136").Value = -33099.999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H50").Value = 12200
$ws.Range("J50").Value = 12200
$ws.Range("L50").Value = 12200
$ws.Range("N50").Value = -13462
$ws.Range("H51").Value = 11000
$ws.Range("I51").Value = 11000
$ws.Range("K51").Value = 11000
$ws.Range("M51").Value = -10490
$ws.Range("H52").Value = 22268.6
$ws.Range("I52").Value = 18135.75
$ws.Range("J52").Value = 38800
$ws.Range("K52").Value = 18135.75
$ws.Range("L52").Value = 38800
$ws.Range("M52").Value = -17909.75
$ws.Range("N52").Value = -39252
$ws.Range("H58").Value = 35966.332
$ws.Range("I58").Value = 6949.5
$ws.Range("K58").Value = 6949.5
$ws.Range("M58").Value = -6641.5
$ws.Range("H122").Value = 3303.0454
$ws.Range("I122").Value = 3129.9473
$ws.Range("K122").Value = 9389.841899999999
$ws.Range("M122").Value = -6939.841899999999
$ws.Range("H132").Value = 4534
$ws.Range("I132").Value = 4409.7334
$ws.Range("K132").Value = 13229.2002
$ws.Range("M132").Value = -10699.2002
